# Fixed Stimulus Absolute Timestamps
# Renames the task-order sheets and updates the stimulus-file timestamps
# referenced in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (new task-order run timestamps) ---------------------
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477900866292"
$wb.Worksheets.Item(2).Name = "NB_TO-16504779024545095"
$wb.Worksheets.Item(3).Name = "RS_TO-16504779024555125"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477902517622"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650477902597622"

# --- Sheet 1: GNG ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504779008374276.csv"
$ws1.Range("B3").Value = "GNG_stims-16504779008492923.csv"
$ws1.Range("B4").Value = "go_stims-1650477900850293.csv"
$ws1.Range("B5").Value = "GNG_stims-16504779008653271.csv"

# --- Sheet 2: NB ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504779015512938.csv"
$ws2.Range("B3").Value = "TB-16504779021853266.csv"
$ws2.Range("B4").Value = "ZB-match_1-16504779010212982.csv"
$ws2.Range("B5").Value = "OB-1650477902050293.csv"
$ws2.Range("B6").Value = "TB-16504779023723285.csv"
$ws2.Range("B7").Value = "TB-16504779024063253.csv"
$ws2.Range("B8").Value = "ZB-match_2-16504779010742977.csv"
$ws2.Range("B9").Value = "ZB-match_4-16504779011262965.csv"
$ws2.Range("B10").Value = "OB-16504779013922956.csv"

# --- Sheet 3: RS ------------------------------------------------------
# (no content changes for this sheet)

# --- Sheet 4: TOL --------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504779024856215.csv"
$ws4.Range("B3").Value = "ZM_stims-16504779024595833.csv"
$ws4.Range("B4").Value = "MM_stims-16504779025016217.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477902486588.csv"
$ws4.Range("B6").Value = "MM_stims-1650477902517622.csv"
$ws4.Range("B7").Value = "ZM_stims-16504779025025935.csv"

# --- Sheet 5: vSAT ---------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504779025255876.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504779025656056.csv"
$ws5.Range("B4").Value = "SAT_stims-16504779025496233.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650477902581625.csv"
